# counter verändert kein direkter reset
# Adds a new "Sheet2" (pasted test-run log/status table) after "Sheet1",
# moves the selection on Sheet1, and makes Sheet2 the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Sheets.Item(1)

# --- Sheet1: just move the selection (no data changes) ---
$ws1.Range("C6").Select()

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$dot = [char]0x00B7
$fw = "pirmicboard_david/unittests/7080gconnect/tmp_main_works_on_gpstracker.py at ccd93bc5e59aa4cae0e1f3a7a5dfe34171ab0504 " + $dot + " Qrist0ph/pirmicboard_david (github.com)"
$hlTarget = "https://github.com/Qrist0ph/pirmicboard_david/blob/ccd93bc5e59aa4cae0e1f3a7a5dfe34171ab0504/unittests/7080gconnect/tmp_main_works_on_gpstracker.py"

# --- Header row ---
$ws2.Range("A1").Value = "Start Bedingungen"
$ws2.Range("B1").Value = "Board"
$ws2.Range("C1").Value = "Firmware"
$ws2.Range("D1").Value = "Sim CCID"
$ws2.Range("E1").Value = "1nce"
$ws2.Range("F1").Value = "Counter"

# --- Firmware hyperlinks first so the "Hyperlink" cell style is created
#     early (keeps the resulting style order close to the source file) ---
$ws2.Hyperlinks.Add($ws2.Range("C2"), $hlTarget, [Type]::Missing, [Type]::Missing, $hlTarget)
$ws2.Range("C2").Value = $fw
$ws2.Hyperlinks.Add($ws2.Range("C3"), $hlTarget, [Type]::Missing, [Type]::Missing, $hlTarget)
$ws2.Range("C3").Value = $fw
$ws2.Hyperlinks.Add($ws2.Range("C4"), $hlTarget, [Type]::Missing, [Type]::Missing, $hlTarget)
$ws2.Range("C4").Value = $fw

# --- CCID / status-code cells: force the leading apostrophe (text-looking
#     number / "+"-prefixed text) just like typing them in manually ---
$ws2.Range("D2").Formula = "'8988228066603839868"
$ws2.Range("D3").Formula = "'8988228066603839867"
$ws2.Range("D4").Formula = "'8988228066603839867"
$ws2.Range("H3").Formula = "'+CPSI: NO SERVICE,Online"
$ws2.Range("H4").Formula = "'+CPSI: LTE NB-IOT,Online,262-01,0xE2A4,35589386,84,EUTRAN-BAND8,3739,0,0,-10,-78,-67,11"

# --- Remaining plain cells ---
$ws2.Range("B2").Value = 2
$ws2.Range("E2").Value = "2023-09-26 07:43:58 UTC"

$ws2.Range("A3").Value = "Sim getauscht"
$ws2.Range("B3").Value = 2
$ws2.Range("E3").Value = "2023-09-26 07:51:53 UTC"
$ws2.Range("F3").Value = 29
$ws2.Range("G3").Value = "No Service, aber scheint dennoch online nach Runde 30"

$ws2.Range("A4").Value = "10 Minutn laufen lassen"
$ws2.Range("B4").Value = 2
$ws2.Range("E4").Value = "2023-09-26 08:06:55 UTC"

# --- Timestamps use the small grey "PortalFont" look ---
$dateRange = $ws2.Range("E2:E4")
$dateRange.Font.Name = "PortalFont"
$dateRange.Font.Size = 10
$dateRange.Font.Color = 4868682

# --- Column widths (best-fit-ish, approximating the pasted-table sizing) ---
$ws2.Columns.Item(1).ColumnWidth = 41.66
$ws2.Columns.Item(4).ColumnWidth = 21.33
$ws2.Columns.Item(5).ColumnWidth = 20.22
$ws2.Columns.Item(7).ColumnWidth = 46.66

# --- Sheet2 becomes the active sheet/tab, with E4 selected ---
$ws2.Range("E4").Select()
